# Commit: "Added WinForm to execute Driver Script"
#
# The PASS/FAIL results that used to be written into the "Results" column
# (column D on the TestCases sheet, column H on the TestSteps sheet) are no
# longer populated from this workbook now that a WinForm drives/executes the
# script, so that stored result data is cleared out. The active sheet/tab
# selection also moves from "TestSteps" back to "TestCases".

$wb = $excel.ActiveWorkbook

$wsTestCases = $wb.Worksheets.Item("TestCases")
$wsTestSteps = $wb.Worksheets.Item("TestSteps")

# Clear out the old PASS/PASS/FAIL results recorded on TestCases!D2:D4
$wsTestCases.Range("D2:D4").ClearContents() | Out-Null

# Clear out the old PASS results recorded on TestSteps!H2:H33
$wsTestSteps.Range("H2:H33").ClearContents() | Out-Null

# TestSteps was the active/selected tab, with G35 selected - update its
# selection to reflect the now-empty H2:H33 block (without leaving it as the
# active tab).
$wsTestSteps.Range("H2:H33").Select() | Out-Null

# TestCases becomes the active tab again, with A2 selected.
$wsTestCases.Activate() | Out-Null
$wsTestCases.Range("A2").Select() | Out-Null
